$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Row 20: 2D_16e_adad_cell1_SNR20dB_50k
$ws.Range("B20").Value = "2D"
$ws.Range("D20").Value = 0.05
$ws.Range("F20").Value = "ad"
$ws.Range("G20").Value = "ad"
$ws.Range("H20").Value = "50k"
$ws.Range("I20").Value = 1
$ws.Range("J20").Value = "[0.05 0.3]"
$ws.Range("L20").Value = "[1 1]"
$ws.Range("K20").Value = "3 layers [2, 1 , 1 ; 2, 1, 0.4; 2,1, 0.2]"
$ws.Range("M20").Value = "no"
$ws.Range("N20").Value = 20

# Row 21: 3D_16e_adad_cell1_SNR20dB_50k
$ws.Range("B21").Value = "3D"
$ws.Range("D21").Value = 0.05
$ws.Range("F21").Value = "ad"
$ws.Range("G21").Value = "ad"
$ws.Range("H21").Value = "50k"
$ws.Range("I21").Value = 1
$ws.Range("J21").Value = "[0.05 0.3]"
$ws.Range("K21").Value = "3 layers [2, 1 , 1 ; 2, 1, 0.4; 2,1, 0.2]"
$ws.Range("L21").Value = "[1 1]"
$ws.Range("M21").Value = "no"
$ws.Range("N21").Value = 20

$excel.Calculate()

$ws.Range("L25").Select()

$wb.Save()
